$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2125748502994012
$ws.Range("C2").Value = 0.5538922155688623
$ws.Range("J2").Value = 0.01796407185628742
$ws.Range("P2").Value = 0.1437125748502994
$ws.Range("S2").Value = 0.0718562874251497
$ws.Range("B3").Value = 0.0155440414507772
$ws.Range("C3").Value = 0.03626943005181347
$ws.Range("J3").Value = 0.0310880829015544
$ws.Range("P3").Value = 0.6943005181347151
$ws.Range("S3").Value = 0.2227979274611399
$ws.Range("J4").Value = 0.09302325581395349
$ws.Range("P4").Value = 0.6046511627906976
$ws.Range("S4").Value = 0.3023255813953488
$ws.Range("B6").Value = 0.09259259259259259
$ws.Range("D6").Value = 0.009259259259259259
$ws.Range("E6").Value = 0.004629629629629629
$ws.Range("F6").Value = 0.06944444444444445
$ws.Range("J6").Value = 0.212962962962963
$ws.Range("O6").Value = 0.02314814814814815
$ws.Range("Q6").Value = 0.1898148148148148
$ws.Range("R6").Value = 0.04166666666666666
$ws.Range("S6").Value = 0.3564814814814815
$ws.Range("B7").Value = 0.13
$ws.Range("D7").Value = 0.005
$ws.Range("F7").Value = 0.035
$ws.Range("J7").Value = 0.165
$ws.Range("Q7").Value = 0.15
$ws.Range("R7").Value = 0.095
$ws.Range("S7").Value = 0.42
$ws.Range("B8").Value = 0.1064814814814815
$ws.Range("D8").Value = 0.03472222222222222
$ws.Range("E8").Value = 0.002314814814814815
$ws.Range("F8").Value = 0.07175925925925926
$ws.Range("J8").Value = 0.1296296296296296
$ws.Range("O8").Value = 0.01851851851851852
$ws.Range("Q8").Value = 0.162037037037037
$ws.Range("R8").Value = 0.08101851851851852
$ws.Range("S8").Value = 0.3935185185185185
$ws.Range("B9").Value = 0.08936170212765958
$ws.Range("F9").Value = 0.08085106382978724
$ws.Range("J9").Value = 0.09361702127659574
$ws.Range("O9").Value = 0.01702127659574468
$ws.Range("Q9").Value = 0.148936170212766
$ws.Range("R9").Value = 0.1148936170212766
$ws.Range("S9").Value = 0.4553191489361702
$ws.Range("B10").Value = 0.1180781758957655
$ws.Range("D10").Value = 0.02035830618892508
$ws.Range("E10").Value = 0.0008143322475570033
$ws.Range("F10").Value = 0.0732899022801303
$ws.Range("J10").Value = 0.1197068403908795
$ws.Range("O10").Value = 0.01547231270358306
$ws.Range("Q10").Value = 0.1775244299674267
$ws.Range("R10").Value = 0.0732899022801303
$ws.Range("S10").Value = 0.4014657980456026
$ws.Range("G11").Value = 0.1301115241635688
$ws.Range("J11").Value = 0.08550185873605948
$ws.Range("K11").Value = 0.2044609665427509
$ws.Range("L11").Value = 0.5762081784386617
$ws.Range("S11").Value = 0.003717472118959108
$ws.Range("G12").Value = 0.84472049689441
$ws.Range("J12").Value = 0.1180124223602484
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.02484472049689441
$ws.Range("S12").Value = 0.006211180124223602
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2448979591836735
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F15").Value = 0.02538071065989848
$ws.Range("H15").Value = 0.1725888324873096
$ws.Range("I15").Value = 0.1116751269035533
$ws.Range("J15").Value = 0.299492385786802
$ws.Range("K15").Value = 0.06091370558375635
$ws.Range("M15").Value = 0.01522842639593909
$ws.Range("O15").Value = 0.05583756345177665
$ws.Range("S15").Value = 0.2588832487309645
$ws.Range("F16").Value = 0.009803921568627451
$ws.Range("H16").Value = 0.1372549019607843
$ws.Range("I16").Value = 0.07352941176470588
$ws.Range("J16").Value = 0.4362745098039216
$ws.Range("K16").Value = 0.107843137254902
$ws.Range("M16").Value = 0.02450980392156863
$ws.Range("N16").Value = 0.004901960784313725
$ws.Range("O16").Value = 0.04411764705882353
$ws.Range("S16").Value = 0.1617647058823529
$ws.Range("F17").Value = 0.0155440414507772
$ws.Range("H17").Value = 0.1917098445595855
$ws.Range("I17").Value = 0.08808290155440414
$ws.Range("J17").Value = 0.4145077720207254
$ws.Range("K17").Value = 0.06735751295336788
$ws.Range("M17").Value = 0.02072538860103627
$ws.Range("N17").Value = 0.002590673575129534
$ws.Range("O17").Value = 0.05699481865284974
$ws.Range("S17").Value = 0.1424870466321244
$ws.Range("F18").Value = 0.005494505494505495
$ws.Range("H18").Value = 0.1153846153846154
$ws.Range("I18").Value = 0.1153846153846154
$ws.Range("J18").Value = 0.4010989010989011
$ws.Range("K18").Value = 0.1813186813186813
$ws.Range("M18").Value = 0.02197802197802198
$ws.Range("O18").Value = 0.03846153846153846
$ws.Range("S18").Value = 0.1208791208791209
$ws.Range("F19").Value = 0.01298701298701299
$ws.Range("H19").Value = 0.2131398013750955
$ws.Range("I19").Value = 0.1084797555385791
$ws.Range("J19").Value = 0.3758594346829641
$ws.Range("K19").Value = 0.09243697478991597
$ws.Range("M19").Value = 0.02291825821237586
$ws.Range("N19").Value = 0.001527883880825057
$ws.Range("O19").Value = 0.06264323911382735
$ws.Range("S19").Value = 0.1100076394194041
